# "4 mdelo melhores rstds" - reorder the model rows by their (presumably newly
# computed) rstd ranking and refresh all metric columns (r2, r2_test, r2_val,
# r2_vt, mse, mse_test, mse_val, mse_vt) with the freshly recomputed values.
#
# The new ranking moves model_10_6_0 to the top, model_10_6_24 to the bottom,
# and reshuffles model_10_6_23 / model_10_6_11 a few places; every row ends up
# sharing the same refreshed metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order (top -> bottom) of the model labels for rows 2..26.
$labels = @(
    "model_10_6_0",
    "model_10_6_22",
    "model_10_6_21",
    "model_10_6_20",
    "model_10_6_19",
    "model_10_6_18",
    "model_10_6_17",
    "model_10_6_16",
    "model_10_6_15",
    "model_10_6_14",
    "model_10_6_13",
    "model_10_6_23",
    "model_10_6_12",
    "model_10_6_10",
    "model_10_6_9",
    "model_10_6_8",
    "model_10_6_7",
    "model_10_6_6",
    "model_10_6_5",
    "model_10_6_4",
    "model_10_6_3",
    "model_10_6_2",
    "model_10_6_1",
    "model_10_6_11",
    "model_10_6_24"
)

# Refreshed metric values (r2, r2_test, r2_val, r2_vt, mse, mse_test, mse_val, mse_vt)
# shared by every row after the recompute.
$vals = @(
    0.6731329884640765,
    -0.3154660929854114,
    0.9853257330522036,
    0.4742282943685033,
    0.3617455065250397,
    1.285408496856689,
    0.02459991723299026,
    0.6920872926712036
)

$firstRow = 2
$rowCount = $labels.Count
$colCount = $vals.Count

$data = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $data[$i, $j] = $vals[$j]
    }
}

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}

$lastRow = $firstRow + $rowCount - 1
$ws.Range("B$firstRow`:I$lastRow").Value = $data
